$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1981.6666
$ws.Range("I40").Value = 1995
$ws.Range("J40").Value = 1975
$ws.Range("K40").Value = 1995
$ws.Range("L40").Value = 1975
$ws.Range("M40").Value = -1820
$ws.Range("N40").Value = -2325
$ws.Range("H64").Value = 3372.8
$ws.Range("I64").Value = 3313.6365
$ws.Range("J64").Value = 3445.111
$ws.Range("K64").Value = 3313.6365
$ws.Range("L64").Value = 3445.111
$ws.Range("M64").Value = -3065.6365
$ws.Range("N64").Value = -3941.111
$ws.Range("H67").Value = 3372.8
$ws.Range("I67").Value = 3313.6365
$ws.Range("J67").Value = 3445.111
$ws.Range("K67").Value = 3313.6365
$ws.Range("L67").Value = 3445.111
$ws.Range("M67").Value = -2455.6365
$ws.Range("N67").Value = -5161.111
$ws.Range("H88").Value = 6000
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -4812
$ws.Range("H91").Value = 6000
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -6808
$ws.Range("H113").Value = 2524.9375
$ws.Range("I113").Value = 2423
$ws.Range("K113").Value = 2423
$ws.Range("M113").Value = 831
$ws.Range("H125").Value = 5500
$ws.Range("I125").Value = 10000
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 90000
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -87540
$ws.Range("N125").Value = -13920
$ws.Range("H138").Value = 2104.2542
$ws.Range("I138").Value = 1253
$ws.Range("J138").Value = 3261.96
$ws.Range("K138").Value = 3759
$ws.Range("L138").Value = 9785.880000000001
$ws.Range("M138").Value = 1381
$ws.Range("N138").Value = -20065.88

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1101.375
$ws.Range("I2").Value = 1081.4
$ws.Range("J2").Value = 1134.6666
$ws.Range("K2").Value = 1081.4
$ws.Range("L2").Value = 1134.6666
$ws.Range("M2").Value = -968.4000000000001
$ws.Range("N2").Value = -1360.6666
$ws.Range("H32").Value = 10316.469
$ws.Range("I32").Value = 9650.075000000001
$ws.Range("J32").Value = 13527.272
$ws.Range("K32").Value = 9650.075000000001
$ws.Range("L32").Value = 13527.272
$ws.Range("M32").Value = -9363.075000000001
$ws.Range("N32").Value = -14101.272
$ws.Range("H97").Value = 4281.75
$ws.Range("I97").Value = 4533.25
$ws.Range("J97").Value = 3527.25
$ws.Range("K97").Value = 4533.25
$ws.Range("L97").Value = 3527.25
$ws.Range("M97").Value = -4037.25
$ws.Range("N97").Value = -4519.25
$ws.Range("H102").Value = 2217.7778
$ws.Range("I102").Value = 2303.3333
$ws.Range("J102").Value = 2175
$ws.Range("K102").Value = 2303.3333
$ws.Range("L102").Value = 2175
$ws.Range("M102").Value = -681.3332999999998
$ws.Range("N102").Value = -5419
$ws.Range("H106").Value = 52370
$ws.Range("J106").Value = 52370
$ws.Range("L106").Value = 52370
$ws.Range("N106").Value = -54894
$ws.Range("H116").Value = 1101.375
$ws.Range("I116").Value = 1081.4
$ws.Range("J116").Value = 1134.6666
$ws.Range("K116").Value = 1081.4
$ws.Range("L116").Value = 1134.6666
$ws.Range("M116").Value = 1212.6
$ws.Range("N116").Value = -5722.6666
$ws.Range("H132").Value = 7578305
$ws.Range("I132").Value = 15627457
$ws.Range("J132").Value = 2632.5881
$ws.Range("K132").Value = 46882371
$ws.Range("L132").Value = 7897.7643
$ws.Range("M132").Value = -46879841
$ws.Range("N132").Value = -12957.7643

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1101.375
$ws.Range("I3").Value = 1081.4
$ws.Range("J3").Value = 1134.6666
$ws.Range("K3").Value = 1081.4
$ws.Range("L3").Value = 1134.6666
$ws.Range("M3").Value = -967.4000000000001
$ws.Range("N3").Value = -1362.6666
$ws.Range("H134").Value = 4954.467
$ws.Range("I134").Value = 4109.8335
$ws.Range("J134").Value = 8333
$ws.Range("K134").Value = 12329.5005
$ws.Range("L134").Value = 24999
$ws.Range("M134").Value = -9794.500499999998
$ws.Range("N134").Value = -30069

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1121.5
$ws.Range("I16").Value = 1070.1
$ws.Range("J16").Value = 1250
$ws.Range("K16").Value = 1070.1
$ws.Range("L16").Value = 1250
$ws.Range("M16").Value = -783.0999999999999
$ws.Range("N16").Value = -1824
$ws.Range("H31").Value = 5652950
$ws.Range("I31").Value = 3473.7546
$ws.Range("J31").Value = 55556656
$ws.Range("K31").Value = 3473.7546
$ws.Range("L31").Value = 55556656
$ws.Range("M31").Value = -3178.7546
$ws.Range("N31").Value = -55557246
$ws.Range("H34").Value = 5652950
$ws.Range("I34").Value = 3473.7546
$ws.Range("J34").Value = 55556656
$ws.Range("K34").Value = 3473.7546
$ws.Range("L34").Value = 55556656
$ws.Range("M34").Value = -3271.7546
$ws.Range("N34").Value = -55557060
$ws.Range("H105").Value = 1544.4445
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 975
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 975
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -4469
$ws.Range("H107").Value = 848.1875
$ws.Range("I107").Value = 571.4
$ws.Range("K107").Value = 571.4
$ws.Range("M107").Value = 1348.6
$ws.Range("H113").Value = 1121.5
$ws.Range("I113").Value = 1070.1
$ws.Range("J113").Value = 1250
$ws.Range("K113").Value = 1070.1
$ws.Range("L113").Value = 1250
$ws.Range("M113").Value = 1099.9
$ws.Range("N113").Value = -5590

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1278.45
$ws.Range("J34").Value = 1542.4375
$ws.Range("L34").Value = 4627.3125
$ws.Range("N34").Value = -4795.3125
$ws.Range("H39").Value = 500.05713
$ws.Range("J39").Value = 500.05713
$ws.Range("L39").Value = 1500.17139
$ws.Range("N39").Value = -2088.17139
$ws.Range("H51").Value = 177.81818
$ws.Range("I51").Value = 177.81818
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 533.4545400000001
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -73.45454000000007
$ws.Range("N51").ClearContents()
$ws.Range("H55").Value = 542.7619
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 542.7619
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 1628.2857
$ws.Range("N55").Value = -1982.2857
$ws.Range("M55").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 781.9394
$ws.Range("J113").Value = 1187.0714
$ws.Range("L113").Value = 3561.2142
$ws.Range("N113").Value = -7901.2142
$ws.Range("H117").Value = 2000
$ws.Range("J117").Value = 2000
$ws.Range("L117").Value = 6000
$ws.Range("N117").Value = -12884
$ws.Range("H125").Value = 3547.1428
$ws.Range("J125").Value = 4033.3333
$ws.Range("L125").Value = 12099.9999
$ws.Range("N125").Value = -21939.9999
$ws.Range("H129").Value = 3522.3713
$ws.Range("J129").Value = 3065.8572
$ws.Range("L129").Value = 9197.571599999999
$ws.Range("N129").Value = -19197.5716
$ws.Range("H132").Value = 760.5263
$ws.Range("J132").Value = 975
$ws.Range("L132").Value = 8775
$ws.Range("N132").Value = -13835

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 201474.2
$ws.Range("I113").Value = 334670.34
$ws.Range("J113").Value = 1680
$ws.Range("K113").Value = 334670.34
$ws.Range("L113").Value = 1680
$ws.Range("M113").Value = -332500.34
$ws.Range("N113").Value = -6020
$ws.Range("H126").Value = 3533.0212
$ws.Range("I126").Value = 2622.8928
$ws.Range("J126").Value = 4874.263
$ws.Range("K126").Value = 7868.678400000001
$ws.Range("L126").Value = 14622.789
$ws.Range("M126").Value = -5398.678400000001
$ws.Range("N126").Value = -19562.789
$ws.Range("H132").Value = 5834.8335
$ws.Range("I132").Value = 5601.1113
$ws.Range("J132").Value = 5975.067
$ws.Range("K132").Value = 16803.3339
$ws.Range("L132").Value = 17925.201
$ws.Range("M132").Value = -14273.3339
$ws.Range("N132").Value = -22985.201

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7157
$ws.Range("I7").Value = 12999.4
$ws.Range("K7").Value = 12999.4
$ws.Range("M7").Value = -12887.4
$ws.Range("H126").Value = 7157
$ws.Range("I126").Value = 12999.4
$ws.Range("K126").Value = 38998.2
$ws.Range("M126").Value = -36528.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1680.5161
$ws.Range("I132").Value = 1201.2
$ws.Range("J132").Value = 2129.875
$ws.Range("K132").Value = 3603.6
$ws.Range("L132").Value = 6389.625
$ws.Range("M132").Value = -1073.6
$ws.Range("N132").Value = -11449.625
